$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '67.141.26'
$ws.Range("E2").Value = '  +3.12%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '3.768.04'
$ws.Range("E3").Value = '  +7.46%  '

# Row 4: TetherUSD
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.03%  '

# Row 5: BNB
$ws.Range("D5").Value = '420.73'
$ws.Range("E5").Value = '  +0.38%  '

# Row 6: Solana
$ws.Range("D6").Value = '131.97'
$ws.Range("E6").Value = '  -0.80%  '

# Row 7: LidoStakedEther
$ws.Range("D7").Value = '3.760.68'
$ws.Range("E7").Value = '  +7.56%  '

# Row 8: XRP
$ws.Range("D8").Value = '0.649'
$ws.Range("E8").Value = '  -0.86%  '

# Row 9: USDC
$ws.Range("E9").Value = '  +0.01%  '

# Row 10: Cardano
$ws.Range("D10").Value = '0.773'
$ws.Range("E10").Value = '  -0.81%  '

# Row 11: Dogecoin
$ws.Range("E11").Value = '  +13.23%  '

# Row 12: ShibaInu
$ws.Range("D12").Value = '0.0000414'
$ws.Range("E12").Value = '  +57.27%  '

# Row 13: Avalanche
$ws.Range("E13").Value = '  -1.67%  '

# Row 14: Polkadot
$ws.Range("D14").Value = '10.39'
$ws.Range("E14").Value = '  +2.97%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = '4.359.70'
$ws.Range("E15").Value = '  +7.34%  '

# Row 16: TRON
$ws.Range("E16").Value = '  -1.08%  '

# Row 17: WrappedEther
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.775.78'
$ws.Range("E17").Value = '  +7.44%  '

# Row 18: Chainlink
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").Value = '20.71'
$ws.Range("E18").Value = '  +0.79%  '

# Row 19: Uniswap
$ws.Range("D19").Value = '13.26'
$ws.Range("E19").Value = '  +3.76%  '

# Row 20: Polygon
$ws.Range("E20").Value = '  +2.89%  '

# Row 21: WrappedBTC
$ws.Range("D21").Value = '67.154.07'
$ws.Range("E21").Value = '  +3.28%  '

# Row 22: BitcoinCash
$ws.Range("D22").Value = '446.58'
$ws.Range("E22").Value = '  -2.03%  '

# Row 23: InternetComputer(DFINITY)
$ws.Range("D23").Value = '16.00'
$ws.Range("E23").Value = '  +19.56%  '

# Row 24: Litecoin
$ws.Range("D24").Value = '89.43'
$ws.Range("E24").Value = '  -0.95%  '

# Row 25: ImmutableX
$ws.Range("D25").Value = '3.13'
$ws.Range("E25").Value = '  -2.81%  '

# Row 26: EthereumClassic
$ws.Range("D26").Value = '38.94'
$ws.Range("E26").Value = '  +14.00%  '

# Row 27: PancakeSwap
$ws.Range("B27").Value = 'PancakeSwap'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D27").Value = '3.34'
$ws.Range("E27").Value = '  -3.17%  '

# Row 28: Filecoin
$ws.Range("B28").Value = 'Filecoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D28").Value = '10.17'
$ws.Range("E28").Value = '  +2.04%  '

# Row 29: LEO
$ws.Range("E29").Value = '  +5.30%  '

# Row 30: Hedera
$ws.Range("E30").Value = '  +6.90%  '

# Row 31: Cosmos
$ws.Range("D31").Value = '12.69'
$ws.Range("E31").Value = '  +0.76%  '

# Row 32: Toncoin
$ws.Range("D32").Value = '2.71'
$ws.Range("E32").Value = '  -1.07%  '

# Row 33: RenderToken
$ws.Range("D33").Value = '7.28'
$ws.Range("E33").Value = '  -3.32%  '

# Row 34: Kaspa
$ws.Range("D34").Value = '0.166'
$ws.Range("E34").Value = '  +1.52%  '

# Row 35: InjectiveProtocol
$ws.Range("D35").Value = '41.96'
$ws.Range("E35").Value = '  +4.46%  '

# Row 36: OKB
$ws.Range("D36").Value = '56.91'
$ws.Range("E36").Value = '  -0.29%  '

# Row 37: Dai
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  -0.02%  '

# Row 38: VeChain
$ws.Range("E38").Value = '  -3.38%  '

# Row 39: PEPE
$ws.Range("D39").Value = '0.0₃0770'
$ws.Range("E39").Value = '  +3.92%  '

# Row 40: Stellar
$ws.Range("E40").Value = '  -0.38%  '

# Row 41: ThetaToken
$ws.Range("D41").Value = '2.98'
$ws.Range("E41").Value = '  +27.54%  '

# Row 42: EnergySwap
$ws.Range("D42").Value = '28.75'
$ws.Range("E42").Value = '  +31.70%  '

# Row 43: FirstDigitalUSD
$ws.Range("D43").Value = '0.997'
$ws.Range("E43").Value = '  -0.17%  '

# Row 44: LidoDAOToken
$ws.Range("D44").Value = '3.46'
$ws.Range("E44").Value = '  +4.51%  '

# Row 45: ARBITRUM
$ws.Range("B45").Value = 'ARBITRUM'
$ws.Range("C45").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D45").Value = '2.17'
$ws.Range("E45").Value = '  +7.76%  '

# Row 46: Monero
$ws.Range("B46").Value = 'Monero'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D46").Value = '147.64'
$ws.Range("E46").Value = '  +1.03%  '

# Row 47: ApeXProtocol
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").Value = '3.17'
$ws.Range("E47").Value = '  +24.30%  '

# Row 48: Stacks
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").Value = '2.90'
$ws.Range("E48").Value = '  -5.73%  '

# Row 49: WEMIXToken
$ws.Range("B49").Value = 'WEMIXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D49").Value = '2.66'
$ws.Range("E49").Value = '  -3.49%  '

# Row 50: NEARProtocol
$ws.Range("D50").Value = '4.35'
$ws.Range("E50").Value = '  -4.77%  '

# Row 51: TheGraph
$ws.Range("E51").Value = '  -2.52%  '
